$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 5 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(5,1).PasteSpecial(-4122)
$ws.Cells.Item(5,1).Value = 45023
$ws.Cells.Item(5,2).Value = 10
$ws.Cells.Item(5,3).Value = 11
$ws.Cells.Item(5,4).Value = "Консультирование по вопросу ликвидации организации"
$ws.Cells.Item(5,5).Value = "токио"

# ---- Row 6 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$ws.Cells.Item(6,1).Value = 45023
$ws.Cells.Item(6,2).Value = 15
$ws.Cells.Item(6,3).Value = 16
$ws.Cells.Item(6,4).Value = "Консультирование по вопросу ликвидации организации"
$ws.Cells.Item(6,5).Value = "волтерс"

# ---- Row 7 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(7,1).PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = 45026
$ws.Cells.Item(7,2).Value = 9
$ws.Cells.Item(7,3).Value = 10
$ws.Cells.Item(7,4).Value = "Консультирование по вопросу ликвидации организации"
$ws.Cells.Item(7,5).Value = "фис"

# ---- Row 8 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(8,1).PasteSpecial(-4122)
$ws.Cells.Item(8,1).Value = 45026
$ws.Cells.Item(8,2).Value = 11
$ws.Cells.Item(8,3).Value = 15
$ws.Cells.Item(8,4).Value = "Консультирование по вопросу ликвидации организации"
$ws.Cells.Item(8,5).Value = "токио"

# ---- Row 9 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(9,1).PasteSpecial(-4122)
$ws.Cells.Item(9,1).Value = 45026
$ws.Cells.Item(9,2).Value = 17
$ws.Cells.Item(9,3).Value = 18
$ws.Cells.Item(9,4).Value = "Консультирование по вопросу ликвидации организации"
$ws.Cells.Item(9,5).Value = "волтерс"

# ---- Row 10 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$ws.Cells.Item(10,1).Value = 45029
$ws.Cells.Item(10,2).Value = 10
$ws.Cells.Item(10,3).Value = 11
$ws.Cells.Item(10,4).Value = "Подготовка договора купли-продажи"
$ws.Cells.Item(10,5).Value = "ЭС"

# ---- Row 11 ----
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(11,1).PasteSpecial(-4122)
$ws.Cells.Item(11,1).Value = 45030
$ws.Cells.Item(11,2).Value = 15
$ws.Cells.Item(11,3).Value = 16
$ws.Cells.Item(11,4).Value = "Поиск информации по отказу от капремонта"
$ws.Cells.Item(11,5).Value = "ЭС"

# ---- Rows 12-16 (blank, formatted like column A date cells, spans 1:5) ----
$ws.Cells.Item(2,1).Copy()
$ws.Range("A12:A16").PasteSpecial(-4122)

# ---- Rows 17-26 (blank, formatted like column A date cells, spans 1:1) ----
$ws.Cells.Item(2,1).Copy()
$ws.Range("A17:A26").PasteSpecial(-4122)

# ---- Update selection ----
$ws.Range("C18").Select()

Write-Output "done"
